$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 2-9 (columns D, L, M, N, O, P, S) after re-sorting
# the weekly price entries by date.
$rows = @(
    @{ Row = 2; D = 44446; L = "Primera";  M = 60; N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 3; D = 44452; L = "Primera";  M = 60; N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 4; D = 44448; L = "Primera";  M = 60; N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 5; D = 44460; L = "Especial"; M = 60; N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 6; D = 44460; L = "Primera";  M = 30; N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 7; D = 44487; L = "Primera";  M = 30; N = 23000; O = 24000; P = 23500; S = 2350 },
    @{ Row = 8; D = 44461; L = "Especial"; M = 60; N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 9; D = 44461; L = "Primera";  M = 30; N = 30000; O = 30000; P = 30000; S = 3000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("S$row").Value = $r.S
}
